# Tag the "quiz answer" shapes on three slides with alt text
# "QuizAnswer" so they can be programmatically identified later
# (e.g. answer-reveal automation). This mirrors the author's
# "set alternative text" action applied identically across the
# affected slides/shapes.

function Find-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Name -eq $name) {
            return $shp
        }
    }
    return $null
}

$p = $ppt.ActivePresentation

# Slide with sldId 598 — smiley-face answer marker.
$s598 = $p.Slides.Item(5)
(Find-ShapeByName $s598 "Smiley Face 5").AlternativeText = "QuizAnswer"

# Slide with sldId 599 — smiley-face marker and the "Table 8" answer table.
$s599 = $p.Slides.Item(12)
(Find-ShapeByName $s599 "Smiley Face 5").AlternativeText = "QuizAnswer"
(Find-ShapeByName $s599 "Table 8").AlternativeText = "QuizAnswer"

# Slide with sldId 602 — smiley-face answer marker.
$s602 = $p.Slides.Item(14)
(Find-ShapeByName $s602 "Smiley Face 4").AlternativeText = "QuizAnswer"
